$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.196163585097679
$ws.Range("C2").Value = 0.3482182018094591
$ws.Range("D2").Value = 0.03380982488847906
$ws.Range("F2").Value = 0.2442865934455405
$ws.Range("G2").Value = 0.1146285456894631
$ws.Range("H2").Value = 0.3002454461834354
$ws.Range("O2").Value = 0.7222467790411145

$ws.Range("B3").Value = 1.045275677438326
$ws.Range("C3").Value = 0.3101611371751005
$ws.Range("D3").Value = 0.02949114967595534
$ws.Range("F3").Value = 0.247929670491402
$ws.Range("G3").Value = 0.1183797293561675
$ws.Range("H3").Value = 0.3061218083232511
$ws.Range("O3").Value = 0.7422668822984875

$ws.Range("B4").Value = 0.952239542118491
$ws.Range("C4").Value = 0.2867206712037387
$ws.Range("D4").Value = 0.02682921527202353
$ws.Range("F4").Value = 0.2505023081252347
$ws.Range("G4").Value = 0.1209357462871559
$ws.Range("H4").Value = 0.3099805673773908
$ws.Range("O4").Value = 0.7556121575679953

$ws.Range("B5").Value = 0.9142313355809506
$ws.Range("C5").Value = 0.2771509102665561
$ws.Range("D5").Value = 0.02574196051907052
$ws.Range("F5").Value = 0.2516348236892583
$ws.Range("G5").Value = 0.1220405337469757
$ws.Range("H5").Value = 0.3116159737164601
$ws.Range("O5").Value = 0.7613142159335098

$ws.Range("B6").Value = 0.9079144371643793
$ws.Range("C6").Value = 0.2755608219278542
$ws.Range("D6").Value = 0.02556127445284062
$ws.Range("F6").Value = 0.2518279506901813
$ws.Range("G6").Value = 0.1222277884136993
$ws.Range("H6").Value = 0.3118913296972501
$ws.Range("O6").Value = 0.7622769389053374

$ws.Range("B7").Value = 0.9517273318160733
$ws.Range("C7").Value = 0.2865916801198978
$ws.Range("D7").Value = 0.02681456216534173
$ws.Range("F7").Value = 0.2505172413099359
$ws.Range("G7").Value = 0.1209503904785763
$ws.Range("H7").Value = 0.3100023683640245
$ws.Range("O7").Value = 0.7556879909310581

$ws.Range("B8").Value = 1.144220156659799
$ws.Range("C8").Value = 0.3351118266829758
$ws.Range("D8").Value = 0.0323229213014713
$ws.Range("F8").Value = 0.2454728638543848
$ws.Range("G8").Value = 0.1158692673996882
$ws.Range("H8").Value = 0.3022195360999405
$ws.Range("O8").Value = 0.7289305953758358

$ws.Range("B9").Value = 1.518490512564824
$ws.Range("C9").Value = 0.4296459678457722
$ws.Range("D9").Value = 0.0430403489289688
$ws.Range("F9").Value = 0.238258071423779
$ws.Range("G9").Value = 0.1079273917558154
$ws.Range("H9").Value = 0.2889503646041476
$ws.Range("O9").Value = 0.6848572886466755

$ws.Range("B10").Value = 1.791390635605353
$ws.Range("C10").Value = 0.4986896639481984
$ws.Range("D10").Value = 0.05085953244424957
$ws.Range("F10").Value = 0.2346074683331523
$ws.Range("G10").Value = 0.1033477703701422
$ws.Range("H10").Value = 0.2804218772328113
$ws.Range("O10").Value = 0.6576555415277312

$ws.Range("B11").Value = 1.915065521416523
$ws.Range("C11").Value = 0.5300031208670362
$ws.Range("D11").Value = 0.05440408326026613
$ws.Range("F11").Value = 0.2333088495291307
$ws.Range("G11").Value = 0.1015417252977926
$ws.Range("H11").Value = 0.2768081839197407
$ws.Range("O11").Value = 0.6464178919169541

$ws.Range("B12").Value = 1.961828051682119
$ws.Range("C12").Value = 0.5418463344281577
$ws.Range("D12").Value = 0.055744451548847
$ws.Range("F12").Value = 0.2328694649366838
$ws.Range("G12").Value = 0.1008980798925094
$ws.Range("H12").Value = 0.2754781185731474
$ws.Range("O12").Value = 0.6423269652474488

$ws.Range("B13").Value = 1.951760086451145
$ws.Range("C13").Value = 0.5392963457313158
$ws.Range("D13").Value = 0.05545586409567704
$ws.Range("F13").Value = 0.2329617600019205
$ws.Range("G13").Value = 0.1010349033753002
$ws.Range("H13").Value = 0.2757628639262109
$ws.Range("O13").Value = 0.6432006848052509

$ws.Range("B14").Value = 1.918914132568148
$ws.Range("C14").Value = 0.5309777649951002
$ws.Range("D14").Value = 0.05451439427018556
$ws.Range("F14").Value = 0.233271649669696
$ws.Range("G14").Value = 0.101487962975682
$ws.Range("H14").Value = 0.2766979891435923
$ws.Range("O14").Value = 0.6460780253196674

$ws.Range("B15").Value = 1.898785770492339
$ws.Range("C15").Value = 0.5258804783799746
$ws.Range("D15").Value = 0.05393746987917325
$ws.Range("F15").Value = 0.2334682958675387
$ws.Range("G15").Value = 0.1017707304418209
$ws.Range("H15").Value = 0.2772757796622187
$ws.Range("O15").Value = 0.6478619380471429

$ws.Range("B16").Value = 1.783298496349573
$ws.Range("C16").Value = 0.4966412722491782
$ws.Range("D16").Value = 0.05062762995108017
$ws.Range("F16").Value = 0.2346996482893218
$ws.Range("G16").Value = 0.10347141257472
$ws.Range("H16").Value = 0.2806634024026948
$ws.Range("O16").Value = 0.6584129123229729

$ws.Range("B17").Value = 1.712328550445875
$ws.Range("C17").Value = 0.4786790458705923
$ws.Range("D17").Value = 0.04859390265431784
$ws.Range("F17").Value = 0.2355479969667513
$ws.Range("G17").Value = 0.1045860189257581
$ws.Range("H17").Value = 0.2828098012961995
$ws.Range("O17").Value = 0.6651774726398401

$ws.Range("B18").Value = 1.671464574797994
$ws.Range("C18").Value = 0.4683387705326254
$ws.Range("D18").Value = 0.04742299129345895
$ws.Range("F18").Value = 0.23607000500418
$ws.Range("G18").Value = 0.1052531833305501
$ws.Range("H18").Value = 0.2840693794549409
$ws.Range("O18").Value = 0.669175182595815

$ws.Range("B19").Value = 1.657621281769138
$ws.Range("C19").Value = 0.4648362328985058
$ws.Range("D19").Value = 0.0470263433914937
$ws.Range("F19").Value = 0.2362525883380897
$ws.Range("G19").Value = 0.1054835395449025
$ws.Range("H19").Value = 0.2845001456778427
$ws.Range("O19").Value = 0.6705470676773331

$ws.Range("B20").Value = 1.719887992203724
$ws.Range("C20").Value = 0.4805920806941231
$ws.Range("D20").Value = 0.04881051763119615
$ws.Range("F20").Value = 0.2354541614283363
$ws.Range("G20").Value = 0.1044646659078659
$ws.Range("H20").Value = 0.2825787224902072
$ws.Range("O20").Value = 0.6644463005451087

$ws.Range("B21").Value = 1.928563722056424
$ws.Range("C21").Value = 0.5334215345798157
$ws.Range("D21").Value = 0.05479097839719316
$ws.Range("F21").Value = 0.2331792036824112
$ws.Range("G21").Value = 0.1013537924680321
$ws.Range("H21").Value = 0.2764222781367991
$ws.Range("O21").Value = 0.6452284067772496

$ws.Range("B22").Value = 2.064533251881301
$ws.Range("C22").Value = 0.5678636459757058
$ws.Range("D22").Value = 0.05868857892804158
$ws.Range("F22").Value = 0.2319978060207504
$ws.Range("G22").Value = 0.09955553238255987
$ws.Range("H22").Value = 0.2726223422665726
$ws.Range("O22").Value = 0.6336278948680985

$ws.Range("B23").Value = 1.992002467529801
$ws.Range("C23").Value = 0.5494893050973246
$ws.Range("D23").Value = 0.05660938946157046
$ws.Range("F23").Value = 0.2326002925092894
$ws.Range("G23").Value = 0.1004936760929382
$ws.Range("H23").Value = 0.2746299363706513
$ws.Range("O23").Value = 0.6397311540709012

$ws.Range("B24").Value = 1.71647056274827
$ws.Range("C24").Value = 0.4797272397983079
$ws.Range("D24").Value = 0.04871259127162375
$ws.Range("F24").Value = 0.2354964777169499
$ws.Range("G24").Value = 0.1045194475216142
$ws.Range("H24").Value = 0.2826831135323076
$ws.Range("O24").Value = 0.6647765250219777

$ws.Range("B25").Value = 1.417596464407666
$ws.Range("C25").Value = 0.404141479218822
$ws.Range("D25").Value = 0.04015039802149545
$ws.Range("F25").Value = 0.2399213783982681
$ws.Range("G25").Value = 0.1098570558486038
$ws.Range("H25").Value = 0.2923261063439817
$ws.Range("O25").Value = 0.6958750777134171
